# Cargo-bay material list update:
#   - "MI_Trim_Metal_Panels_A" to-do note now says "[to do in unreal] ..."
#   - "Interior Metal Panel Trim" note gets a trailing qualifier
#   - leftover proofing marks (spell/grammar check artifacts) are cleared
#     as part of the rewrite, same as Word does when content is retyped.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "[to do] MI_Trim_Metal_Panels_A", $true, $false, $false, $false, $false,
    $true, 1, $false, "[to do in unreal] MI_Trim_Metal_Panels_A", 2) | Out-Null

$d.Content.Find.Execute(
    "Interior Metal Panel Trim", $true, $false, $false, $false, $false,
    $true, 1, $false, "Interior Metal Panel Trim (the one from a couple semesters ago)", 2) | Out-Null

# Rebuild the whole body from clean paragraph XML so no stray
# <w:proofErr/> spell/grammar-check markers survive the edit.
$paras = @(
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "Material names across DCCs:" },
    @{ Style = $null;           Ilvl = $null; Text = "" },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "MI_Trim_Metal_Ship_Large_A " },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "Used for the exterior of ship main painted metal panels" },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "MI_Trim_Metal_Ship_Large_B " },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "Used for the mechanical parts" },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "[to do in unreal] MI_Trim_Metal_Panels_A" },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "Interior Metal Panel Trim (the one from a couple semesters ago)" },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "MI_Ship_Glass " },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "Used for ship glass" },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "MI_Ship_Decals" },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "Labels for ship." },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "[to do]MI_Trim_Metal_Interior" },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "A sharper version used for the interior" },
    @{ Style = "ListParagraph"; Ilvl = "0"; Text = "X MI_Trim_Concrete_Outdoor_Large" },
    @{ Style = "ListParagraph"; Ilvl = "1"; Text = "A" }
)

$body = ""
foreach ($para in $paras) {
    if ($null -eq $para.Style) {
        $body += "<w:p/>"
    } else {
        $escaped = $para.Text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $body += '<w:p><w:pPr><w:pStyle w:val="' + $para.Style + '"/><w:numPr><w:ilvl w:val="' + $para.Ilvl + '"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p>'
    }
}

$xmlFrag = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document>'

$d.Content.InsertXML($xmlFrag)
